# Tyotunnit.xlsx edit: add a new work-session row (row 30: 2020-09-21)
# and refresh the active selection / scroll position to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# --- Fill in the new entry (date, start time, end time) on row 30 ---
# The existing D/F/G formulas on this row are shared formulas that were
# already present (referencing an empty A/B/C) - entering real data here
# just makes them recalculate, which Excel's autorecalc handles for us.

$ws.Range("A30").Value2 = 44095
$ws.Range("A30").NumberFormat = "d-mmm"

$ws.Range("B30").Value2 = 0.64583333333333337
$ws.Range("B30").NumberFormat = "h:mm"

$ws.Range("C30").Value2 = 0.89583333333333337
$ws.Range("C30").NumberFormat = "h:mm"

# --- Update the view state to where the author ended up editing ---
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("E26").Select() | Out-Null
